# Auto update Excel log 2026-02-04 14:18:06
# Appends new sensor-log rows to the PIR, Humidity, and Temperature sheets,
# mirroring continued data capture from the logging device.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($sheetName, $rows)

    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($r in $rows) {
        $rowNum = $r[0]

        # Force the whole row to Text format first so date-/time-/percent-
        # looking strings are stored verbatim instead of being reinterpreted
        # as dates, times, or numbers by Excel's auto-detection.
        $ws.Range("A" + $rowNum + ":F" + $rowNum).NumberFormat = "@"

        $ws.Cells.Item($rowNum, 1).Value = $r[1]
        $ws.Cells.Item($rowNum, 2).Value = $r[2]
        $ws.Cells.Item($rowNum, 3).Value = $r[3]
        $ws.Cells.Item($rowNum, 4).Value = $r[4]
        $ws.Cells.Item($rowNum, 5).Value = $r[5]
        $ws.Cells.Item($rowNum, 6).Value = $r[6]
    }
}

# PIR sheet: rows 192-203 (Date, Timestamp, Hour, Location, Value, Status)
$pirRows = @(
    @(192,"2026-02-04","14:17:02","14:00","Bathroom","No Motion","Inactive"),
    @(193,"2026-02-04","14:17:04","14:00","Bathroom","Motion Detected","Active"),
    @(194,"2026-02-04","14:17:11","14:00","Bathroom","No Motion","Inactive"),
    @(195,"2026-02-04","14:17:15","14:00","Bathroom","Motion Detected","Active"),
    @(196,"2026-02-04","14:17:24","14:00","Bathroom","No Motion","Inactive"),
    @(197,"2026-02-04","14:17:29","14:00","Bathroom","Motion Detected","Active"),
    @(198,"2026-02-04","14:17:35","14:00","Bathroom","No Motion","Inactive"),
    @(199,"2026-02-04","14:17:40","14:00","Bathroom","No Motion","Inactive"),
    @(200,"2026-02-04","14:17:45","14:00","Bathroom","No Motion","Inactive"),
    @(201,"2026-02-04","14:17:50","14:00","Bathroom","No Motion","Inactive"),
    @(202,"2026-02-04","14:17:55","14:00","Bathroom","No Motion","Inactive"),
    @(203,"2026-02-04","14:18:00","14:00","Bathroom","No Motion","Inactive")
)
Add-LogRows "PIR" $pirRows

# Humidity sheet: rows 158-170
$humidityRows = @(
    @(158,"2026-02-04","14:17:01","14:00","Bathroom","77.8%","Active"),
    @(159,"2026-02-04","14:17:03","14:00","Bathroom","76.8%","Active"),
    @(160,"2026-02-04","14:17:07","14:00","Bathroom","77.7%","Active"),
    @(161,"2026-02-04","14:17:12","14:00","Bathroom","76.8%","Active"),
    @(162,"2026-02-04","14:17:17","14:00","Bathroom","77.8%","Active"),
    @(163,"2026-02-04","14:17:22","14:00","Bathroom","76.8%","Active"),
    @(164,"2026-02-04","14:17:27","14:00","Bathroom","77.7%","Active"),
    @(165,"2026-02-04","14:17:32","14:00","Bathroom","76.9%","Active"),
    @(166,"2026-02-04","14:17:37","14:00","Bathroom","77.9%","Active"),
    @(167,"2026-02-04","14:17:42","14:00","Bathroom","77.0%","Active"),
    @(168,"2026-02-04","14:17:47","14:00","Bathroom","77.9%","Active"),
    @(169,"2026-02-04","14:17:52","14:00","Bathroom","76.8%","Active"),
    @(170,"2026-02-04","14:17:57","14:00","Bathroom","77.8%","Active")
)
Add-LogRows "Humidity" $humidityRows

# Temperature sheet: rows 158-170
$temperatureRows = @(
    @(158,"2026-02-04","14:17:02","14:00","Bathroom","24.8C","Active"),
    @(159,"2026-02-04","14:17:03","14:00","Bathroom","24.7C","Active"),
    @(160,"2026-02-04","14:17:08","14:00","Bathroom","24.7C","Active"),
    @(161,"2026-02-04","14:17:13","14:00","Bathroom","24.7C","Active"),
    @(162,"2026-02-04","14:17:18","14:00","Bathroom","24.7C","Active"),
    @(163,"2026-02-04","14:17:23","14:00","Bathroom","24.7C","Active"),
    @(164,"2026-02-04","14:17:28","14:00","Bathroom","24.7C","Active"),
    @(165,"2026-02-04","14:17:33","14:00","Bathroom","24.7C","Active"),
    @(166,"2026-02-04","14:17:38","14:00","Bathroom","24.7C","Active"),
    @(167,"2026-02-04","14:17:43","14:00","Bathroom","24.6C","Active"),
    @(168,"2026-02-04","14:17:48","14:00","Bathroom","24.6C","Active"),
    @(169,"2026-02-04","14:17:53","14:00","Bathroom","24.6C","Active"),
    @(170,"2026-02-04","14:17:58","14:00","Bathroom","24.7C","Active")
)
Add-LogRows "Temperature" $temperatureRows
